$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (2-7) and write new data rows (2-8)
$ws.Range("A2:E7").ClearContents()

# Row 2 (Ticket Type left blank on purpose)
$ws.Range("B2").Value = "Story test One"
$ws.Range("C2").Value = "Sample Story Ticket"
$ws.Range("D2").Value = "Junior Erick Otieno"
$ws.Range("E2").Value = "Medium"

# Row 3
$ws.Range("A3").Value = "Task"
$ws.Range("B3").Value = "Task Test One"
$ws.Range("C3").Value = "Sample Task ticket"
$ws.Range("D3").Value = "Junior Erick Otieno"
$ws.Range("E3").Value = "Minor"

# Row 4
$ws.Range("A4").Value = "Bug"
$ws.Range("C4").Value = "No Bug summary"
$ws.Range("D4").Value = "Anthony Njuguna"
$ws.Range("E4").Value = "Blocker"

# Row 5
$ws.Range("A5").Value = "Story"
$ws.Range("B5").Value = "Client Walk through feature"
$ws.Range("C5").Value = "Create a walkthrough pop up to guide new clients"
$ws.Range("D5").Value = "Pesh Kirigo"
$ws.Range("E5").Value = "Must Have"

# Row 6
$ws.Range("A6").Value = "Bug"
$ws.Range("B6").Value = "Register button not working"
$ws.Range("D6").Value = "Junior Erick Otieno"
$ws.Range("E6").Value = "Critical"

# Row 7
$ws.Range("A7").Value = "Task"
$ws.Range("B7").Value = "Dark theme"
$ws.Range("C7").Value = "Include dark theme for the site"
$ws.Range("E7").Value = "Should Have"

# Row 8 (new row)
$ws.Range("A8").Value = "Task"
$ws.Range("B8").Value = "Include Remember Me checkbox"
$ws.Range("C8").Value = "Include remember me checkbox to allow customers to be automatically logged in during subsequent logins"
$ws.Range("D8").Value = "Junior Erick Otieno"
$ws.Range("E8").Value = "Would Like"

# Wrap text + row height for row 8
$ws.Range("C8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 29

# Selection
$ws.Range("E8").Select()

# Page setup
$ws.PageSetup.Orientation = 1
